# Leveling_base + crit chance 1% => 10%
#
# Insert a new "nothing" (base/no-weapon) entry at the top of the guns
# leveling table on the "guns" sheet, pushing every existing weapon row
# down by one. Values are written explicitly (rather than via a native
# row-insert/shift) to keep the untouched numeric literals byte-identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("guns")

# Shift existing data rows (old row 2..12) down to (new row 3..13),
# working from the bottom up so nothing is overwritten before it's moved.
$ws.Range("A13").Value = "Witch Killer's Blade"
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 1.25

$ws.Range("A12").Value = "nightingale's dagger"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 1.35

$ws.Range("A11").Value = "meteorite sword"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 1.25

$ws.Range("A10").Value = "war axe"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 1.25

$ws.Range("A9").Value = "axe"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 1.1499999999999999

$ws.Range("A8").Value = "mace"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 1.05

$ws.Range("A7").Value = "steel sword"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 1.1000000000000001

$ws.Range("A6").Value = "iron sword"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 1.07

$ws.Range("A5").Value = "copper sword"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1.03

$ws.Range("A4").Value = "dagger"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1.1499999999999999

$ws.Range("A3").Value = "knife"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1.05

# New base row: no weapon equipped, no level requirement, no crit bonus.
$ws.Range("A2").Value = "nothing"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# Match the saved view state after the edit.
$ws.Activate()
$ws.Range("D2").Select()
